$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -20.292
    6   = -22.291
    7   = -19.851
    8   = -22.146
    16  = -21.93
    20  = -20.137
    21  = -20.057
    28  = -21.893
    29  = -21.344
    30  = -21.91699999999999
    32  = -21.648
    40  = -20.174
    46  = -21.87
    51  = -21.98
    52  = -22.084
    57  = -22.263
    59  = -22.395
    62  = -22.263
    66  = -21.616
    73  = -20.082
    74  = -21.112
    77  = -20.342
    92  = -21.50600000000001
    100 = -22.217
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
